# GJES-321 - Update upload formats in Create Guest List screen.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Note " -> "Note/Affiliation"
$ws.Range("D1").Value = "Note/Affiliation"

# "Addtional guests " -> "Number of addtional guests  "
$ws.Range("E1").Value = "Number of addtional guests  "

# Drop the "F for Female/ M for Male" header entirely - column F has no label now
$ws.Range("F1").ClearContents()

# Widen column E so the longer header text fits
$ws.Columns.Item(5).ColumnWidth = 29
